# Add "2022-Q4" sheet data and update the "总计" (summary) sheet
# to reflect the newly added quarter, per commit "feat: add 2022-Q4 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet (sheet 1): insert a 2022-Q4 row at
#    the top of the data (row 2), push the existing rows down by one,
#    and grow the table by one row (now 4 data rows instead of 3).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Extend the A-column "index" style down into the new row 5 by copying
# the existing formatted cell (keeps border/alignment consistent).
$summary.Range("A4").Copy($summary.Range("A5"))

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.02

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 0.58

$summary.Range("B4").Value = "2021-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.89

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2020-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 2.14

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计" and
#    before "2022-Q3" (so the tab order becomes 总计, 2022-Q4, 2022-Q3,
#    2021-Q1, 2020-Q4), then fill it with the quarter's fund holdings.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)          # "2022-Q3" (insertion point)
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Copy the column layout/styles from the (now shifted) "2022-Q3" sheet so
# the header row, borders and row-index column match the rest of the
# workbook exactly, then overwrite every cell with this quarter's data.
$q3 = $wb.Worksheets.Item(3)                   # "2022-Q3" after the insert
$q3.Range("A1:H6").Copy($q4.Range("A1"))

# Fund-code-like and percentage-like columns are stored as text in the
# source data (e.g. "008778", "0.60"), so force text format before
# writing them to avoid Excel auto-converting to numbers.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$q4.Range("B2").Value = "008778"
$q4.Range("C2").Value = "嘉实中证500指数增强A"
$q4.Range("D2").Value = "0.60"
$q4.Range("E2").Value = "93.52"
$q4.Range("F2").Value = "1.67"
$q4.Range("G2").Value = "0.0100"
$q4.Range("H2").Value = 10

$q4.Range("B3").Value = "008779"
$q4.Range("C3").Value = "嘉实中证500指数增强C"
$q4.Range("D3").Value = "0.40"
$q4.Range("E3").Value = "93.52"
$q4.Range("F3").Value = "1.67"
$q4.Range("G3").Value = "0.0067"
$q4.Range("H3").Value = 10

$q4.Range("B4").Value = "851088"
$q4.Range("C4").Value = "海通量化成长精选一年持有期混合A"
$q4.Range("D4").Value = "0.39"
$q4.Range("E4").Value = "83.15"
$q4.Range("F4").Value = "0.80"
$q4.Range("G4").Value = "0.0031"
$q4.Range("H4").Value = 10

$q4.Range("B5").Value = "850010"
$q4.Range("C5").Value = "海通量化成长精选一年持有期混合B"
$q4.Range("D5").Value = "0.26"
$q4.Range("E5").Value = "83.15"
$q4.Range("F5").Value = "0.80"
$q4.Range("G5").Value = "0.0021"
$q4.Range("H5").Value = 10

$q4.Range("B6").Value = "851099"
$q4.Range("C6").Value = "海通量化成长精选一年持有期混合C"
$q4.Range("D6").Value = "0.03"
$q4.Range("E6").Value = "83.15"
$q4.Range("F6").Value = "0.80"
$q4.Range("G6").Value = "0.0002"
$q4.Range("H6").Value = 10
